$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "Matteo pilati"
$ws.Range("B68").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C68").Value = "Michele Merighi | Clitoriders"
$ws.Range("D68").Value = "Matteo Diener | U.SGUARNA"
$ws.Range("E68").Value = "Moris Benedetti | Gli Introvabili"
$ws.Range("F68").Value = "Emiliano Bici | Power Ginger"
